# Update NATMI ligand-receptor pair stats (Cd34-Selp) with new TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 227.282303
$ws.Range("H2").Value = 681.846909
$ws.Range("I2").Value = 0.6094595465130797
$ws.Range("J2").Value = 0.6094595465130795
$ws.Range("M2").Value = 35.585194
$ws.Range("N2").Value = 106.755582
$ws.Range("O2").Value = 0.9972091466993565
$ws.Range("P2").Value = 0.9972091466993567
$ws.Range("Q2").Value = 8087.884845021782
$ws.Range("R2").Value = 72790.96360519604
$ws.Range("S2").Value = 0.6077586343260849
$ws.Range("T2").Value = 0.6077586343260849

# Row 3
$ws.Range("G3").Value = 227.282303
$ws.Range("H3").Value = 681.846909
$ws.Range("I3").Value = 0.6094595465130797
$ws.Range("J3").Value = 0.6094595465130795
$ws.Range("M3").Value = 0.093901
$ws.Range("N3").Value = 0.281703
$ws.Range("O3").Value = 0.002631401590341653
$ws.Range("P3").Value = 0.002631401590341654
$ws.Range("Q3").Value = 21.342035534003
$ws.Range("R3").Value = 192.078319806027
$ws.Range("S3").Value = 0.001603732819943421
$ws.Range("T3").Value = 0.001603732819943421

# Row 4
$ws.Range("G4").Value = 227.282303
$ws.Range("H4").Value = 681.846909
$ws.Range("I4").Value = 0.6094595465130797
$ws.Range("J4").Value = 0.6094595465130795
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.00569
$ws.Range("N4").Value = 0.01707
$ws.Range("O4").Value = 0.0001594517103017434
$ws.Range("P4").Value = 0.0001594517103017434
$ws.Range("Q4").Value = 1.29323630407
$ws.Range("R4").Value = 11.63912673663
$ws.Range("S4").Value = 0.00009717936705123549
$ws.Range("T4").Value = 0.00009717936705123549

# Row 5
$ws.Range("I5").Value = 0.3727053955221387
$ws.Range("J5").Value = 0.3727053955221385
$ws.Range("M5").Value = 35.585194
$ws.Range("N5").Value = 106.755582
$ws.Range("O5").Value = 0.9972091466993565
$ws.Range("P5").Value = 0.9972091466993567
$ws.Range("Q5").Value = 4946.018710097705
$ws.Range("R5").Value = 44514.16839087934
$ws.Range("S5").Value = 0.371665229438878
$ws.Range("T5").Value = 0.371665229438878

# Row 6
$ws.Range("I6").Value = 0.3727053955221387
$ws.Range("J6").Value = 0.3727053955221385
$ws.Range("M6").Value = 0.093901
$ws.Range("N6").Value = 0.281703
$ws.Range("O6").Value = 0.002631401590341653
$ws.Range("P6").Value = 0.002631401590341654
$ws.Range("Q6").Value = 13.051386003316
$ws.Range("R6").Value = 117.462474029844
$ws.Range("S6").Value = 0.0009807375705058706
$ws.Range("T6").Value = 0.0009807375705058704

# Row 7
$ws.Range("I7").Value = 0.3727053955221387
$ws.Range("J7").Value = 0.3727053955221385
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.00569
$ws.Range("N7").Value = 0.01707
$ws.Range("O7").Value = 0.0001594517103017434
$ws.Range("P7").Value = 0.0001594517103017434
$ws.Range("Q7").Value = 0.79085831204
$ws.Range("R7").Value = 7.117724808359999
$ws.Range("S7").Value = 0.00005942851275469275
$ws.Range("T7").Value = 0.00005942851275469274

# Row 8
$ws.Range("G8").Value = 6.377814666666666
$ws.Range("H8").Value = 19.133444
$ws.Range("I8").Value = 0.0171021675827138
$ws.Range("J8").Value = 0.01710216758271379
$ws.Range("M8").Value = 35.585194
$ws.Range("N8").Value = 106.755582
$ws.Range("O8").Value = 0.9972091466993565
$ws.Range("P8").Value = 0.9972091466993567
$ws.Range("Q8").Value = 226.9557722093786
$ws.Range("R8").Value = 2042.601949884408
$ws.Range("S8").Value = 0.01705443794186742
$ws.Range("T8").Value = 0.01705443794186742

# Row 9
$ws.Range("G9").Value = 6.377814666666666
$ws.Range("H9").Value = 19.133444
$ws.Range("I9").Value = 0.0171021675827138
$ws.Range("J9").Value = 0.01710216758271379
$ws.Range("M9").Value = 0.093901
$ws.Range("N9").Value = 0.281703
$ws.Range("O9").Value = 0.002631401590341653
$ws.Range("P9").Value = 0.002631401590341654
$ws.Range("Q9").Value = 0.5988831750146666
$ws.Range("R9").Value = 5.389948575131998
$ws.Range("S9").Value = 0.00004500267097544256
$ws.Range("T9").Value = 0.00004500267097544255

# Row 10
$ws.Range("G10").Value = 6.377814666666666
$ws.Range("H10").Value = 19.133444
$ws.Range("I10").Value = 0.0171021675827138
$ws.Range("J10").Value = 0.01710216758271379
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.00569
$ws.Range("N10").Value = 0.01707
$ws.Range("O10").Value = 0.0001594517103017434
$ws.Range("P10").Value = 0.0001594517103017434
$ws.Range("Q10").Value = 0.03628976545333332
$ws.Range("R10").Value = 0.3266078890799999
$ws.Range("S10").Value = 0.000002726969870930748
$ws.Range("T10").Value = 0.000002726969870930747

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.2733126666666666
$ws.Range("H11").Value = 0.819938
$ws.Range("I11").Value = 0.0007328903820679218
$ws.Range("J11").Value = 0.0007328903820679217
$ws.Range("M11").Value = 35.585194
$ws.Range("N11").Value = 106.755582
$ws.Range("O11").Value = 0.9972091466993565
$ws.Range("P11").Value = 0.9972091466993567
$ws.Range("Q11").Value = 9.725884265990667
$ws.Range("R11").Value = 87.53295839391599
$ws.Range("S11").Value = 0.0007308449925261176
$ws.Range("T11").Value = 0.0007308449925261177

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.2733126666666666
$ws.Range("H12").Value = 0.819938
$ws.Range("I12").Value = 0.0007328903820679218
$ws.Range("J12").Value = 0.0007328903820679217
$ws.Range("M12").Value = 0.093901
$ws.Range("N12").Value = 0.281703
$ws.Range("O12").Value = 0.002631401590341653
$ws.Range("P12").Value = 0.002631401590341654
$ws.Range("Q12").Value = 0.02566433271266666
$ws.Range("R12").Value = 0.230978994414
$ws.Range("S12").Value = 0.000001928528916919631
$ws.Range("T12").Value = 0.000001928528916919631

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.2733126666666666
$ws.Range("H13").Value = 0.819938
$ws.Range("I13").Value = 0.0007328903820679218
$ws.Range("J13").Value = 0.0007328903820679217
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.00569
$ws.Range("N13").Value = 0.01707
$ws.Range("O13").Value = 0.0001594517103017434
$ws.Range("P13").Value = 0.0001594517103017434
$ws.Range("Q13").Value = 0.001555149073333333
$ws.Range("R13").Value = 0.01399634166
$ws.Range("S13").Value = 0.0000001168606248844283
$ws.Range("T13").Value = 0.0000001168606248844283

